# Users seed workbook update
# - Drop the "phone" column entirely (email, phone, service_user, firm_ref -> email, service_user, firm_ref)
# - Merge the two firm-admin rows into a single "admin@stopoverx.com" service-user row
# - Blank out the now-unused second row (previously the glidequest admin)
# - Drop the per-row firm_ref values (column stays as a header only)
# - Keep customer / service-user / mikhail rows intact

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove every hyperlink up front; we'll re-create the ones that are still
# needed once the sheet has its final shape (deleting any single cell's
# hyperlink collection clears all of them on this engine).
$ws.Range("A1").Hyperlinks.Delete()

# Drop column B ("phone"); service_user and firm_ref shift left.
$ws.Columns.Item(2).Delete()

# Row 2: was the wildwheeladventures admin -> becomes the single stopoverx admin.
$ws.Range("A2").Value = "admin@stopoverx.com"
$ws.Range("B2").Value = $true
$ws.Range("C2").ClearContents()

# Row 3: was the glidequesttours admin -> now blank placeholder row. Keep its
# existing (hyperlink-flavoured) cell style as-is; it's still used below as a
# formatting donor for the untouched rows further down.
$ws.Range("A3").ClearContents()
$ws.Range("B3").Value = $false
$ws.Range("C3").ClearContents()

# Rows 4-10 keep their emails/service_user flags; just drop firm_ref values.
$ws.Range("C4").ClearContents()
$ws.Range("C5").ClearContents()
$ws.Range("C6").ClearContents()
$ws.Range("C7").ClearContents()
$ws.Range("C8").ClearContents()
$ws.Range("C9").ClearContents()
$ws.Range("C10").ClearContents()

# Re-create mailto hyperlinks for every row that still has an email address.
# Adding a hyperlink resets that cell's style to this engine's plain
# "Hyperlink" look, which is exactly what row 2 (a brand-new link) should
# get, so it's left alone. The other rows already had a custom hyperlink
# style before this edit, so restore it afterwards by copying the format
# from A3, whose style was never touched by Hyperlinks.Add.
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:" + $ws.Range("A2").Value()) | Out-Null

$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:" + $ws.Range("A4").Value()) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A5"), "mailto:" + $ws.Range("A5").Value()) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A6"), "mailto:" + $ws.Range("A6").Value()) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:" + $ws.Range("A7").Value()) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A8"), "mailto:" + $ws.Range("A8").Value()) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A9"), "mailto:" + $ws.Range("A9").Value()) | Out-Null
$ws.Hyperlinks.Add($ws.Range("A10"), "mailto:" + $ws.Range("A10").Value()) | Out-Null

$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4").PasteSpecial(-4122) | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null
$ws.Range("A7").PasteSpecial(-4122) | Out-Null
$ws.Range("A8").PasteSpecial(-4122) | Out-Null
$ws.Range("A9").PasteSpecial(-4122) | Out-Null
$ws.Range("A10").PasteSpecial(-4122) | Out-Null

# Put the active selection on A2, matching the saved view state.
$ws.Range("A2").Select() | Out-Null
